$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (BTC)
$ws.Range("B2").Value = 'BTC'
$ws.Range("C2").Value = 'Bitcoin'
$ws.Range("D2").Value = 68983
$ws.Range("E2").Value = 1361579868007
$ws.Range("F2").Value = 17610962012
$ws.Range("G2").Value = -0.37545

# Row 3 (ETH)
$ws.Range("B3").Value = 'ETH'
$ws.Range("C3").Value = 'Ethereum'
$ws.Range("D3").Value = 3800.82
$ws.Range("E3").Value = 457408957998
$ws.Range("F3").Value = 9290941563
$ws.Range("G3").Value = 1.63758

# Row 4 (USDT)
$ws.Range("B4").Value = 'USDT'
$ws.Range("C4").Value = 'Tether'
$ws.Range("D4").Value = 0.9988669999999999
$ws.Range("E4").Value = 111888382394
$ws.Range("F4").Value = 26703735996
$ws.Range("G4").Value = -0.17004

# Row 5 (BNB)
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'BNB'
$ws.Range("D5").Value = 600.49
$ws.Range("E5").Value = 92527029871
$ws.Range("F5").Value = 408008931
$ws.Range("G5").Value = -0.53846

# Row 6 (SOL)
$ws.Range("B6").Value = 'SOL'
$ws.Range("C6").Value = 'Solana'
$ws.Range("D6").Value = 163.15
$ws.Range("E6").Value = 73466971974
$ws.Range("F6").Value = 2054274930
$ws.Range("G6").Value = -3.17117

# Row 7 (STETH)
$ws.Range("B7").Value = 'STETH'
$ws.Range("C7").Value = 'Lido Staked Ether'
$ws.Range("D7").Value = 3799.86
$ws.Range("E7").Value = 35354597227
$ws.Range("F7").Value = 62162050
$ws.Range("G7").Value = 1.72976

# Row 8 (USDC)
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'USDC'
$ws.Range("D8").Value = 0.998881
$ws.Range("E8").Value = 32702003742
$ws.Range("F8").Value = 2933690722
$ws.Range("G8").Value = -0.25052

# Row 9 (XRP)
$ws.Range("B9").Value = 'XRP'
$ws.Range("C9").Value = 'XRP'
$ws.Range("D9").Value = 0.534792
$ws.Range("E9").Value = 29698592087
$ws.Range("F9").Value = 442062348
$ws.Range("G9").Value = -0.51252

# Row 10 (DOGE)
$ws.Range("B10").Value = 'DOGE'
$ws.Range("C10").Value = 'Dogecoin'
$ws.Range("D10").Value = 0.169542
$ws.Range("E10").Value = 24596121062
$ws.Range("F10").Value = 1318436125
$ws.Range("G10").Value = 1.31917

# Row 11 (TON)
$ws.Range("B11").Value = 'TON'
$ws.Range("C11").Value = 'Toncoin'
$ws.Range("D11").Value = 6.3
$ws.Range("E11").Value = 21917181150
$ws.Range("F11").Value = 136791697
$ws.Range("G11").Value = -1.72868

# Row 12 (ADA)
$ws.Range("B12").Value = 'ADA'
$ws.Range("C12").Value = 'Cardano'
$ws.Range("D12").Value = 0.460149
$ws.Range("E12").Value = 16316761736
$ws.Range("F12").Value = 203245408
$ws.Range("G12").Value = -0.28653

# Row 13 (AVAX)
$ws.Range("B13").Value = 'AVAX'
$ws.Range("C13").Value = 'Avalanche'
$ws.Range("D13").Value = 37.1
$ws.Range("E13").Value = 14631138744
$ws.Range("F13").Value = 223073681
$ws.Range("G13").Value = -3.11585

# Row 14 (SHIB)
$ws.Range("B14").Value = 'SHIB'
$ws.Range("C14").Value = 'Shiba Inu'
$ws.Range("D14").Value = 0.00002452
$ws.Range("E14").Value = 14512897960
$ws.Range("F14").Value = 312350066
$ws.Range("G14").Value = -1.43333

# Row 15 (WBTC)
$ws.Range("B15").Value = 'WBTC'
$ws.Range("C15").Value = 'Wrapped Bitcoin'
$ws.Range("D15").Value = 69094
$ws.Range("E15").Value = 10742943518
$ws.Range("F15").Value = 141945337
$ws.Range("G15").Value = -0.15439

# Row 16 (DOT)
$ws.Range("B16").Value = 'DOT'
$ws.Range("C16").Value = 'Polkadot'
$ws.Range("D16").Value = 7.42
$ws.Range("E16").Value = 10177998960
$ws.Range("F16").Value = 141874089
$ws.Range("G16").Value = 1.58825

# Row 17 (LINK)
$ws.Range("B17").Value = 'LINK'
$ws.Range("C17").Value = 'Chainlink'
$ws.Range("D17").Value = 17.23
$ws.Range("E17").Value = 10166094733
$ws.Range("F17").Value = 467212608
$ws.Range("G17").Value = 0.83655

# Row 18 (TRX)
$ws.Range("B18").Value = 'TRX'
$ws.Range("C18").Value = 'TRON'
$ws.Range("D18").Value = 0.113734
$ws.Range("E18").Value = 9957112125
$ws.Range("F18").Value = 233286062
$ws.Range("G18").Value = -0.1324

# Row 19 (BCH)
$ws.Range("B19").Value = 'BCH'
$ws.Range("C19").Value = 'Bitcoin Cash'
$ws.Range("D19").Value = 486.48
$ws.Range("E19").Value = 9621422042
$ws.Range("F19").Value = 160681790
$ws.Range("G19").Value = -1.54505

# Row 20 (UNI)
$ws.Range("B20").Value = 'UNI'
$ws.Range("C20").Value = 'Uniswap'
$ws.Range("D20").Value = 11.49
$ws.Range("E20").Value = 8680989209
$ws.Range("F20").Value = 516609249
$ws.Range("G20").Value = 5.70289

# Row 21 (NEAR)
$ws.Range("B21").Value = 'NEAR'
$ws.Range("C21").Value = 'NEAR Protocol'
$ws.Range("D21").Value = 7.99
$ws.Range("E21").Value = 8634095261
$ws.Range("F21").Value = 474344548
$ws.Range("G21").Value = -0.31056

# Row 22 (PEPE)
$ws.Range("B22").Value = 'PEPE'
$ws.Range("C22").Value = 'Pepe'
$ws.Range("D22").Value = 0.00001603
$ws.Range("E22").Value = 6812573418
$ws.Range("F22").Value = 1216539355
$ws.Range("G22").Value = 4.36425

# Row 23 (MATIC)
$ws.Range("B23").Value = 'MATIC'
$ws.Range("C23").Value = 'Polygon'
$ws.Range("D23").Value = 0.718925
$ws.Range("E23").Value = 6688859440
$ws.Range("F23").Value = 176692953
$ws.Range("G23").Value = -1.17463

# Row 24 (LTC)
$ws.Range("B24").Value = 'LTC'
$ws.Range("C24").Value = 'Litecoin'
$ws.Range("D24").Value = 84.59999999999999
$ws.Range("E24").Value = 6318083207
$ws.Range("F24").Value = 200457470
$ws.Range("G24").Value = -0.37638

# Row 25 (FET)
$ws.Range("B25").Value = 'FET'
$ws.Range("C25").Value = 'Fetch.ai'
$ws.Range("D25").Value = 2.25
$ws.Range("E25").Value = 5696376149
$ws.Range("F25").Value = 132653721
$ws.Range("G25").Value = -3.41661

# Row 26 (ICP)
$ws.Range("B26").Value = 'ICP'
$ws.Range("C26").Value = 'Internet Computer'
$ws.Range("D26").Value = 12.2
$ws.Range("E26").Value = 5682965485
$ws.Range("F26").Value = 57970608
$ws.Range("G26").Value = -1.30919

# Row 27 (LEO)
$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'LEO Token'
$ws.Range("D27").Value = 5.92
$ws.Range("E27").Value = 5495652741
$ws.Range("F27").Value = 1189955
$ws.Range("G27").Value = -0.92383

# Row 28 (DAI)
$ws.Range("B28").Value = 'DAI'
$ws.Range("C28").Value = 'Dai'
$ws.Range("D28").Value = 0.99834
$ws.Range("E28").Value = 5287429821
$ws.Range("F28").Value = 165153793
$ws.Range("G28").Value = -0.17315

# Row 29 (WEETH)
$ws.Range("B29").Value = 'WEETH'
$ws.Range("C29").Value = 'Wrapped eETH'
$ws.Range("D29").Value = 3949.18
$ws.Range("E29").Value = 4727446826
$ws.Range("F29").Value = 37022681
$ws.Range("G29").Value = 1.64397

# Row 30 (ETC)
$ws.Range("B30").Value = 'ETC'
$ws.Range("C30").Value = 'Ethereum Classic'
$ws.Range("D30").Value = 31.67
$ws.Range("E30").Value = 4682969733
$ws.Range("F30").Value = 167394549
$ws.Range("G30").Value = -0.46069

# Row 31 (APT)
$ws.Range("B31").Value = 'APT'
$ws.Range("C31").Value = 'Aptos'
$ws.Range("D31").Value = 8.98
$ws.Range("E31").Value = 3922489425
$ws.Range("F31").Value = 110137057
$ws.Range("G31").Value = -0.95219

# Row 32 (RNDR)
$ws.Range("B32").Value = 'RNDR'
$ws.Range("C32").Value = 'Render'
$ws.Range("D32").Value = 10
$ws.Range("E32").Value = 3894457800
$ws.Range("F32").Value = 98452880
$ws.Range("G32").Value = -2.20133

# Row 33 (HBAR)
$ws.Range("B33").Value = 'HBAR'
$ws.Range("C33").Value = 'Hedera'
$ws.Range("D33").Value = 0.10659
$ws.Range("E33").Value = 3822283349
$ws.Range("F33").Value = 33264651
$ws.Range("G33").Value = -1.95571

# Row 34 (EZETH)
$ws.Range("B34").Value = 'EZETH'
$ws.Range("C34").Value = 'Renzo Restaked ETH'
$ws.Range("D34").Value = 3748.42
$ws.Range("E34").Value = 3741592973
$ws.Range("F34").Value = 64824196
$ws.Range("G34").Value = 2.17674

# Row 35 (IMX)
$ws.Range("B35").Value = 'IMX'
$ws.Range("C35").Value = 'Immutable'
$ws.Range("D35").Value = 2.38
$ws.Range("E35").Value = 3542629757
$ws.Range("F35").Value = 40979232
$ws.Range("G35").Value = -5.1917

# Row 36 (KAS)
$ws.Range("B36").Value = 'KAS'
$ws.Range("C36").Value = 'Kaspa'
$ws.Range("D36").Value = 0.140425
$ws.Range("E36").Value = 3350870056
$ws.Range("F36").Value = 29662531
$ws.Range("G36").Value = 4.92099

# Row 37 (MNT)
$ws.Range("B37").Value = 'MNT'
$ws.Range("C37").Value = 'Mantle'
$ws.Range("D37").Value = 1.021
$ws.Range("E37").Value = 3326088718
$ws.Range("F37").Value = 43819573
$ws.Range("G37").Value = 0.54165

# Row 38 (FIL)
$ws.Range("B38").Value = 'FIL'
$ws.Range("C38").Value = 'Filecoin'
$ws.Range("D38").Value = 5.87
$ws.Range("E38").Value = 3283213723
$ws.Range("F38").Value = 120290762
$ws.Range("G38").Value = -0.20129

# Row 39 (ATOM)
$ws.Range("B39").Value = 'ATOM'
$ws.Range("C39").Value = 'Cosmos Hub'
$ws.Range("D39").Value = 8.359999999999999
$ws.Range("E39").Value = 3272815131
$ws.Range("F39").Value = 146080703
$ws.Range("G39").Value = -1.45652

# Row 40 (CRO)
$ws.Range("B40").Value = 'CRO'
$ws.Range("C40").Value = 'Cronos'
$ws.Range("D40").Value = 0.119629
$ws.Range("E40").Value = 3208782836
$ws.Range("F40").Value = 7181770
$ws.Range("G40").Value = -1.63523

# Row 41 (ARB)
$ws.Range("B41").Value = 'ARB'
$ws.Range("C41").Value = 'Arbitrum'
$ws.Range("D41").Value = 1.2
$ws.Range("E41").Value = 3205113286
$ws.Range("F41").Value = 340733390
$ws.Range("G41").Value = 0.7222499999999999

# Row 42 (XLM)
$ws.Range("B42").Value = 'XLM'
$ws.Range("C42").Value = 'Stellar'
$ws.Range("D42").Value = 0.109276
$ws.Range("E42").Value = 3173165623
$ws.Range("F42").Value = 33229650
$ws.Range("G42").Value = -1.73276

# Row 43 (FDUSD)
$ws.Range("B43").Value = 'FDUSD'
$ws.Range("C43").Value = 'First Digital USD'
$ws.Range("D43").Value = 0.997129
$ws.Range("E43").Value = 3135416280
$ws.Range("F43").Value = 3135682430
$ws.Range("G43").Value = -0.19632

# Row 44 (WIF)
$ws.Range("B44").Value = 'WIF'
$ws.Range("C44").Value = 'dogwifhat'
$ws.Range("D44").Value = 3.04
$ws.Range("E44").Value = 3045225847
$ws.Range("F44").Value = 520778744
$ws.Range("G44").Value = 0.60827

# Row 45 (GRT)
$ws.Range("B45").Value = 'GRT'
$ws.Range("C45").Value = 'The Graph'
$ws.Range("D45").Value = 0.318213
$ws.Range("E45").Value = 3033338471
$ws.Range("F45").Value = 64902188
$ws.Range("G45").Value = -1.91148

# Row 46 (TAO)
$ws.Range("B46").Value = 'TAO'
$ws.Range("C46").Value = 'Bittensor'
$ws.Range("D46").Value = 427.2
$ws.Range("E46").Value = 2925955991
$ws.Range("F46").Value = 17042479
$ws.Range("G46").Value = -1.54337

# Row 47 (OKB)
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'OKB'
$ws.Range("D47").Value = 48.48
$ws.Range("E47").Value = 2912000088
$ws.Range("F47").Value = 4055804
$ws.Range("G47").Value = -0.2425

# Row 48 (STX)
$ws.Range("B48").Value = 'STX'
$ws.Range("C48").Value = 'Stacks'
$ws.Range("D48").Value = 1.98
$ws.Range("E48").Value = 2908053573
$ws.Range("F48").Value = 37016535
$ws.Range("G48").Value = -0.15902

# Row 49 (OP)
$ws.Range("B49").Value = 'OP'
$ws.Range("C49").Value = 'Optimism'
$ws.Range("D49").Value = 2.56
$ws.Range("E49").Value = 2793058470
$ws.Range("F49").Value = 198632404
$ws.Range("G49").Value = 0.79262

# Row 50 (USDE)
$ws.Range("B50").Value = 'USDE'
$ws.Range("C50").Value = 'Ethena USDe'
$ws.Range("D50").Value = 0.999457
$ws.Range("E50").Value = 2780849112
$ws.Range("F50").Value = 71380328
$ws.Range("G50").Value = -0.17412

# Row 51 (MKR)
$ws.Range("B51").Value = 'MKR'
$ws.Range("C51").Value = 'Maker'
$ws.Range("D51").Value = 2822.16
$ws.Range("E51").Value = 2626688601
$ws.Range("F51").Value = 55645838
$ws.Range("G51").Value = 1.27263
